# Generate Report for Handoff
# Update the "Latest Handoff Datetime" (column D) for the b520fca5... row (row 6)
# on both the zh-cn and de-de localization-status sheets, reflecting a fresh
# handoff that has just been generated.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D6").Value = "2016-03-08 10:26:00"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D6").Value = "2016-03-08 10:26:05"
